$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 9.579440671982837 * [Math]::Pow(10, -11)
$ws.Range("E3").Value = 9.579440671982837 * [Math]::Pow(10, -11)

$ws.Range("D4").Value = 0.002166236941271408
$ws.Range("E4").Value = 0.002166236941271408

$ws.Range("D5").Value = 2.94254022519651 * [Math]::Pow(10, -7)
$ws.Range("E5").Value = 2.94254022519651 * [Math]::Pow(10, -7)

$ws.Range("D6").Value = 3.436786180856466 * [Math]::Pow(10, -42)
$ws.Range("E6").Value = 3.436786180856466 * [Math]::Pow(10, -42)

$ws.Range("D7").Value = 0.9999999997692886
$ws.Range("E7").Value = 2.307114499444651 * [Math]::Pow(10, -10)

$ws.Range("D8").Value = 3.061003580600538 * [Math]::Pow(10, -12)
$ws.Range("E8").Value = 0.999999999996939

$ws.Range("D9").Value = 7.27035580250089 * [Math]::Pow(10, -49) * [Math]::Pow(10, 40)
$ws.Range("E9").Value = 0.9999999927296442

$ws.Range("D10").Value = 6.190673795807732 * [Math]::Pow(10, -7)
$ws.Range("E10").Value = 0.9999993809326204

$ws.Range("D11").Value = 6.144811118196574 * [Math]::Pow(10, -13)
$ws.Range("E11").Value = 0.9999999999993855
$ws.Range("F11").Value = 12.66297912597656
